$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.473.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.997.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.05%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "

$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5106"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.49%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4140"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08708"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.132"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.996.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.573"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.430"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.35%  "

$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06514"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("E20").Value = "  +4.53%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.165"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.519.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.95%  "

$ws.Range("E25").Value = "  +1.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.224.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.399"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.11%  "

$ws.Range("E31").Value = "  +0.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1051"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.070"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.853"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.333"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02516"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.427"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.69%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06600"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2204"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.019"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6616"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.233"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6163"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.198"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.664"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.268"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "80.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.60%  "

$ws.Range("E51").Value = "  +1.58%  "

